$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values in Price/Volume/Hora columns
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '245.29'
$ws.Cells.Item(2, 5).Value = '-0.74%'
$ws.Cells.Item(2, 7).Value = '13'

$ws.Cells.Item(3, 4).Value = '29.95'
$ws.Cells.Item(3, 5).Value = '0.16%'
$ws.Cells.Item(3, 7).Value = '13'

$ws.Cells.Item(4, 4).Value = '5.163'
$ws.Cells.Item(4, 5).Value = '-0.14%'
$ws.Cells.Item(4, 7).Value = '13'

$ws.Cells.Item(5, 4).Value = '0.05743'
$ws.Cells.Item(5, 7).Value = '13'

$ws.Cells.Item(6, 5).Value = '0.80%'
$ws.Cells.Item(6, 7).Value = '13'

$ws.Cells.Item(7, 4).Value = '3.291'
$ws.Cells.Item(7, 5).Value = '7.93%'
$ws.Cells.Item(7, 7).Value = '13'

$ws.Cells.Item(8, 4).Value = '0.8576'
$ws.Cells.Item(8, 5).Value = '-0.25%'
$ws.Cells.Item(8, 7).Value = '13'

$ws.Cells.Item(9, 4).Value = '0.8518'
$ws.Cells.Item(9, 5).Value = '-2.10%'
$ws.Cells.Item(9, 7).Value = '13'

$ws.Cells.Item(10, 5).Value = '1.26%'
$ws.Cells.Item(10, 7).Value = '13'

$ws.Cells.Item(11, 4).Value = '0.07098'
$ws.Cells.Item(11, 5).Value = '0.08%'
$ws.Cells.Item(11, 7).Value = '13'

$ws.Cells.Item(12, 4).Value = '0.03145'
$ws.Cells.Item(12, 5).Value = '9.86%'
$ws.Cells.Item(12, 7).Value = '13'

$ws.Cells.Item(13, 4).Value = '0.09372'
$ws.Cells.Item(13, 5).Value = '-0.09%'
$ws.Cells.Item(13, 7).Value = '13'

$ws.Cells.Item(14, 4).Value = '0.001534'
$ws.Cells.Item(14, 5).Value = '1.27%'
$ws.Cells.Item(14, 7).Value = '13'

$ws.Cells.Item(15, 4).Value = '0.0005988'
$ws.Cells.Item(15, 5).Value = '-0.02%'
$ws.Cells.Item(15, 7).Value = '13'

$ws.Cells.Item(16, 4).Value = '0.005960'
$ws.Cells.Item(16, 5).Value = '-3.52%'
$ws.Cells.Item(16, 7).Value = '13'

$ws.Cells.Item(17, 4).Value = '3.537'
$ws.Cells.Item(17, 5).Value = '1.67%'
$ws.Cells.Item(17, 7).Value = '13'

$ws.Cells.Item(18, 4).Value = '2.195'
$ws.Cells.Item(18, 5).Value = '0.68%'
$ws.Cells.Item(18, 7).Value = '13'

$ws.Cells.Item(19, 4).Value = '0.3148'
$ws.Cells.Item(19, 5).Value = '0.08%'
$ws.Cells.Item(19, 7).Value = '13'

$ws.Cells.Item(20, 4).Value = '0.03318'
$ws.Cells.Item(20, 5).Value = '1.89%'
$ws.Cells.Item(20, 7).Value = '13'

$ws.Cells.Item(21, 4).Value = '0.1279'
$ws.Cells.Item(21, 5).Value = '-2.21%'
$ws.Cells.Item(21, 7).Value = '13'

$ws.Cells.Item(22, 4).Value = '3.506'
$ws.Cells.Item(22, 5).Value = '11.34%'
$ws.Cells.Item(22, 7).Value = '13'

$ws.Cells.Item(23, 5).Value = '2.15%'
$ws.Cells.Item(23, 7).Value = '13'

$ws.Cells.Item(24, 5).Value = '-0.48%'
$ws.Cells.Item(24, 7).Value = '13'

$ws.Cells.Item(25, 4).Value = '0.001225'
$ws.Cells.Item(25, 5).Value = '0.33%'
$ws.Cells.Item(25, 7).Value = '13'

$ws.Cells.Item(26, 4).Value = '0.004161'
$ws.Cells.Item(26, 5).Value = '-18.55%'
$ws.Cells.Item(26, 7).Value = '13'

$ws.Cells.Item(27, 5).Value = '-0.88%'
$ws.Cells.Item(27, 7).Value = '13'

$ws.Cells.Item(28, 5).Value = '-25.25%'
$ws.Cells.Item(28, 7).Value = '13'

$ws.Cells.Item(29, 7).Value = '13'

$ws.Cells.Item(30, 7).Value = '13'

$ws.Cells.Item(31, 7).Value = '13'

$ws.Cells.Item(32, 7).Value = '13'

$ws.Cells.Item(33, 7).Value = '13'

$ws.Cells.Item(34, 7).Value = '13'

$ws.Cells.Item(35, 7).Value = '13'

$ws.Cells.Item(36, 7).Value = '13'

$ws.Cells.Item(37, 7).Value = '13'

$ws.Cells.Item(38, 7).Value = '13'

$ws.Cells.Item(39, 7).Value = '13'

$ws.Cells.Item(40, 4).Value = '0.03746'
$ws.Cells.Item(40, 5).Value = '-0.85%'
$ws.Cells.Item(40, 7).Value = '13'

$ws.Cells.Item(41, 4).Value = '0.1072'
$ws.Cells.Item(41, 5).Value = '0.02%'
$ws.Cells.Item(41, 7).Value = '13'

$ws.Cells.Item(42, 2).Value = 'KickToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Cells.Item(42, 4).Value = '0.003582'
$ws.Cells.Item(42, 5).Value = '-37.61%'
$ws.Cells.Item(42, 7).Value = '13'

$ws.Cells.Item(43, 2).Value = 'CEJI'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Cells.Item(43, 4).Value = '0.002459'
$ws.Cells.Item(43, 5).Value = '-5.44%'
$ws.Cells.Item(43, 7).Value = '13'

$ws.Cells.Item(44, 4).Value = '0.009443'
$ws.Cells.Item(44, 5).Value = '-3.64%'
$ws.Cells.Item(44, 7).Value = '13'

$ws.Cells.Item(45, 4).Value = '0.00005475'
$ws.Cells.Item(45, 5).Value = '7.57%'
$ws.Cells.Item(45, 7).Value = '13'

$ws.Cells.Item(46, 5).Value = '-0.01%'
$ws.Cells.Item(46, 7).Value = '13'

$ws.Cells.Item(47, 5).Value = '19.84%'
$ws.Cells.Item(47, 7).Value = '13'

$ws.Cells.Item(48, 4).Value = '0.002219'
$ws.Cells.Item(48, 5).Value = '-19.56%'
$ws.Cells.Item(48, 7).Value = '13'

$ws.Cells.Item(49, 5).Value = '-0.01%'
$ws.Cells.Item(49, 7).Value = '13'

$ws.Cells.Item(50, 5).Value = '-0.01%'
$ws.Cells.Item(50, 7).Value = '13'

$ws.Cells.Item(51, 7).Value = '13'

Write-Output "Applied all changes"